$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (ALC)
$ws.Range("I19").Value = 1054.5714
$ws.Range("L19").Value = 879
$ws.Range("M19").Value = -879.5714
$ws.Range("J19").Value = 879
$ws.Range("K19").Value = 1054.5714
$ws.Range("H19").Value = 973.53845
$ws.Range("N19").Value = -1229

# Row 106 (ALC)
$ws.Range("K106").Value = 4316.6665
$ws.Range("I106").Value = 4316.6665
$ws.Range("H106").Value = 4316.6665
$ws.Range("M106").Value = -3685.6665

# Row 132 (ALC)
$ws.Range("K132").Value = 55653.75
$ws.Range("M132").Value = -53123.75
$ws.Range("H132").Value = 18551.25
$ws.Range("I132").Value = 18551.25

# Row 133 (ALC)
$ws.Range("H133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -160120
$ws.Range("J133").Value = 150000

$ws = $wb.Worksheets.Item("ARM")
# Row 110 (ARM)
$ws.Range("K110").Value = 711.8333
$ws.Range("M110").Value = 1333.1667
$ws.Range("J110").Value = 633
$ws.Range("H110").Value = 696.06665
$ws.Range("L110").Value = 633
$ws.Range("I110").Value = 711.8333
$ws.Range("N110").Value = -4723

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM)
$ws.Range("H94").Value = 1651.4667
$ws.Range("N94").Value = -2601.5
$ws.Range("L94").Value = 1699.5
$ws.Range("J94").Value = 1699.5

# Row 105 (BSM)
$ws.Range("H105").Value = 2777.6
$ws.Range("K105").Value = 2610.5625
$ws.Range("I105").Value = 2610.5625
$ws.Range("M105").Value = -863.5625

# Row 119 (BSM)
$ws.Range("J119").Value = 60000
$ws.Range("N119").Value = -69676
$ws.Range("H119").Value = 60000
$ws.Range("L119").Value = 60000

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (CRP)
$ws.Range("M7").Value = -6216.1763
$ws.Range("I7").Value = 6329.1763
$ws.Range("K7").Value = 6329.1763
$ws.Range("H7").Value = 3920.3225

# Row 31 (CRP)
$ws.Range("L31").Value = 8825.944
$ws.Range("N31").Value = -9415.944
$ws.Range("I31").Value = 1730.5714
$ws.Range("J31").Value = 8825.944
$ws.Range("H31").Value = 6839.24
$ws.Range("K31").Value = 1730.5714
$ws.Range("M31").Value = -1435.5714

# Row 34 (CRP)
$ws.Range("I34").Value = 1730.5714
$ws.Range("L34").Value = 8825.944
$ws.Range("J34").Value = 8825.944
$ws.Range("M34").Value = -1528.5714
$ws.Range("K34").Value = 1730.5714
$ws.Range("H34").Value = 6839.24
$ws.Range("N34").Value = -9229.944

# Row 132 (CRP)
$ws.Range("K132").Value = 13097.6661
$ws.Range("M132").Value = -10567.6661
$ws.Range("H132").Value = 4592.357
$ws.Range("I132").Value = 4365.8887

$ws = $wb.Worksheets.Item("CUL")
# Row 50 (CUL)
$ws.Range("K50").Value = 50639.50199999999
$ws.Range("I50").Value = 16879.834
$ws.Range("M50").Value = -50158.50199999999
$ws.Range("H50").Value = 16879.834

# Row 53 (CUL)
$ws.Range("M53").Value = -50158.50199999999
$ws.Range("H53").Value = 16879.834
$ws.Range("I53").Value = 16879.834
$ws.Range("K53").Value = 50639.50199999999

# Row 70 (CUL)
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73 (CUL)
$ws.Range("H73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("N73").ClearContents()

# Row 87 (CUL)
$ws.Range("K87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("M87").ClearContents()

# Row 90 (CUL)
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("M90").ClearContents()

# Row 116 (CUL)
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

# Row 136 (CUL)
$ws.Range("M136").Value = -8397
$ws.Range("I136").Value = 4499
$ws.Range("K136").Value = 13497
$ws.Range("H136").Value = 7327.75

# Row 138 (CUL)
$ws.Range("H138").Value = 5811.875
$ws.Range("M138").Value = -3856.25
$ws.Range("K138").Value = 8996.25
$ws.Range("I138").Value = 2998.75

$ws = $wb.Worksheets.Item("GSM")
# Row 36 (GSM)
$ws.Range("H36").Value = 733.3333
$ws.Range("K36").Value = 733.3333
$ws.Range("M36").Value = -248.3333
$ws.Range("I36").Value = 733.3333
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# Row 43 (GSM)
$ws.Range("L43").Value = 19800
$ws.Range("J43").Value = 19800
$ws.Range("N43").Value = -20102
$ws.Range("H43").Value = 6200

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (LTW)
$ws.Range("I2").Value = 25971.5
$ws.Range("K2").Value = 25971.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -25859.5
$ws.Range("H2").Value = 25971.5
$ws.Range("N2").ClearContents()

# Row 7 (LTW)
$ws.Range("M7").Value = -6093.6665
$ws.Range("I7").Value = 6205.6665
$ws.Range("K7").Value = 6205.6665
$ws.Range("H7").Value = 7499.25

# Row 16 (LTW)
$ws.Range("K16").Value = 1400
$ws.Range("M16").Value = -1230
$ws.Range("H16").Value = 1400
$ws.Range("I16").Value = 1400

# Row 40 (LTW)
$ws.Range("I40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# Row 61 (LTW)
$ws.Range("K61").Value = 1308.0714
$ws.Range("L61").Value = 7999.75
$ws.Range("H61").Value = 2795.111
$ws.Range("J61").Value = 7999.75
$ws.Range("I61").Value = 1308.0714
$ws.Range("N61").Value = -8403.75
$ws.Range("M61").Value = -1106.0714

# Row 92 (LTW)
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992

# Row 113 (LTW)
$ws.Range("J113").Value = 7999.75
$ws.Range("K113").Value = 1308.0714
$ws.Range("M113").Value = 861.9286
$ws.Range("H113").Value = 2795.111
$ws.Range("I113").Value = 1308.0714
$ws.Range("L113").Value = 7999.75
$ws.Range("N113").Value = -12339.75

# Row 126 (LTW)
$ws.Range("I126").Value = 6205.6665
$ws.Range("H126").Value = 7499.25
$ws.Range("K126").Value = 18616.9995
$ws.Range("M126").Value = -16146.9995

# Row 136 (LTW)
$ws.Range("L136").Value = 14850
$ws.Range("M136").Value = -7293
$ws.Range("I136").Value = 3281
$ws.Range("N136").Value = -19950
$ws.Range("K136").Value = 9843
$ws.Range("J136").Value = 4950
$ws.Range("H136").Value = 4564.846

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (WVR)
$ws.Range("K100").Value = 2696.3332
$ws.Range("I100").Value = 1348.1666
$ws.Range("M100").Value = -2155.3332
$ws.Range("H100").Value = 1441.2858

# Row 122 (WVR)
$ws.Range("L122").Value = 17550
$ws.Range("J122").Value = 5850
$ws.Range("H122").Value = 5381.3335
$ws.Range("N122").Value = -22450

# Row 126 (WVR)
$ws.Range("I126").Value = 4749.5
$ws.Range("H126").Value = 6992.077
$ws.Range("N126").Value = -28906.334
$ws.Range("K126").Value = 14248.5
$ws.Range("M126").Value = -11778.5
$ws.Range("J126").Value = 7988.778
$ws.Range("L126").Value = 23966.334

# Row 136 (WVR)
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -3856.059
$ws.Range("I136").Value = 2135.353
$ws.Range("N136").Value = -19350
$ws.Range("K136").Value = 6406.059
$ws.Range("J136").Value = 4750
$ws.Range("H136").Value = 2972.04

Write-Output "applied changes"